$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (id=2, "Baker, Salas and Nelson Incubator") entirely.
# This shifts all subsequent rows up by one, matching the diff.
$ws.Rows.Item(2).Delete()
